$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 12 data (2017-06-04 -> serial 42890)
$ws.Range("A12").Value = 42890
$ws.Range("B12").Value = 5
$ws.Range("C12").Value = "Reservanto první nasazení, Nahrazení head, keywords, title, Odstraněny nepotřebné stránky, Kontrola správnosti odkazů, Footer a jeho nahrazení ve stránkách, Nasazeno PHP - jednotné menu, Styly, tránka služeb, galerie, 404"

# Match the formatting used by the row above (reuse its date/number styles)
$ws.Range("A11:C11").Copy()
$ws.Range("A12:C12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$wb.Save()
